$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 178, pushing the existing rows 178-279 down to 179-280.
$ws.Rows(178).Insert()

# Populate the new row 178 with the new price-sheet entry.
$ws.Range("A178").Value = 4
$ws.Range("B178").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C178").Value = "Los Lagos"
$ws.Range("D178").Value = 44777
$ws.Range("E178").Value = 10
$ws.Range("F178").Value = "Fruta"
$ws.Range("G178").Value = 100108
$ws.Range("H178").Value = "Tropicales y subtropicales"
$ws.Range("I178").Value = 100108005
$ws.Range("J178").Value = "Piña"
$ws.Range("K178").Value = "Caramelo"
$ws.Range("L178").Value = "Primera"
$ws.Range("M178").Value = 80
$ws.Range("N178").Value = 23000
$ws.Range("O178").Value = 23000
$ws.Range("P178").Value = 23000
$ws.Range("Q178").Value = "$/caja 12 unidades"
$ws.Range("R178").Value = "Ecuador"
$ws.Range("S178").Value = 1917
$ws.Range("T178").Value = 12
